# Nowe testy + Aplikacje konsolowe
# Adds three new timesheet rows (25-27) mirroring the existing row layout
# (date in col A, hours in col B, comment in col C), which extends the
# used range from A1:E24 to A1:E27 and bumps the E2 SUM(B:B) total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: 2020-07-15 (serial 44027), 2.5h
$ws.Range("A25").Value = 44027
$ws.Range("B25").Value = 2.5
$ws.Range("C25").Value = "Praca nad poprawnymi aplikacjami konsolowymi - dla dynamicznego"

# Row 26: 2020-07-16 (serial 44028), 1h
$ws.Range("A26").Value = 44028
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "Uczenie się o testach."

# Row 27: 2020-07-16 (serial 44028), 1h
$ws.Range("A27").Value = 44028
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "Nowe testy"

# Match the saved selection/active cell (C24 -> C27) and scroll position
# (topLeftCell A7 -> A13) recorded in the workbook view.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()
